$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Training Dashboard")
$ws2 = $wb.Worksheets.Item("Exam Dashboard")

# --- Header / title formatting: consolidate bold font to white text ---
# (the two near-duplicate bold fonts -- title "bold 14pt" and header "bold" --
# become a single bold white font shared by the title and the column headers)
$ws1.Range("A1").Font.Size = 11
$ws1.Range("A1").Font.Color = 16777215
$ws1.Range("A2:K2").Font.Color = 16777215

$ws2.Range("A1").Font.Size = 11
$ws2.Range("A1").Font.Color = 16777215
$ws2.Range("A2:G2").Font.Color = 16777215

# --- Training Dashboard: refresh "PERIOD TO EXPIRE" (H) and "LAST UPDATE" (I) ---
$periods = @{
    3  = 675
    4  = 368
    5  = 484
    6  = 490
    7  = 357
    8  = 483
    9  = 355
    10 = 370
    11 = 489
    12 = 697
    13 = 490
    14 = 377
    15 = 490
    16 = 489
    17 = 489
    18 = 489
    19 = 47
    20 = 298
    21 = 298
    22 = 300
    23 = 204
    24 = 205
    25 = 311
    26 = 304
    27 = 311
    28 = 312
    29 = 310
    30 = 332
    31 = 332
}

foreach ($row in $periods.Keys) {
    $ws1.Cells.Item($row, 8).Value = $periods[$row]
}

$lastUpdateRange = $ws1.Range("I3:I31")
$lastUpdateRange.NumberFormat = "@"
$lastUpdateRange.Value = "16-Sep-2025"

# --- Exam Dashboard: update comments and shrink the COMMENTS column ---
$ws2.Range("E3").Value = "date is valid"
$ws2.Range("E4").Value = "date is valid"
$ws2.Range("E5").Value = "date is valid"
$ws2.Range("E6").Value = "date is valid"
$ws2.Range("E7").Value = "date is valid"

$ws2.Columns.Item(5).ColumnWidth = 14.1
